$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1031
$ws1.Range("F7").Value = 1035
$ws1.Range("F11").Value = 613
$ws1.Range("F12").Value = 46
$ws1.Range("F15").Value = 1984
$ws1.Range("G18").Value = "已售罄"
$ws1.Range("F21").Value = 632
$ws1.Range("F26").Value = 3475
$ws1.Range("F32").Value = 486
$ws1.Range("F36").Value = 228
$ws1.Range("F37").Value = 319
$ws1.Range("F38").Value = 818

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 578
$ws2.Range("F6").Value = 314

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1031
$ws4.Range("F9").Value = 578
$ws4.Range("F10").Value = 1035
$ws4.Range("F13").Value = 314
$ws4.Range("F16").Value = 613
$ws4.Range("F17").Value = 46
$ws4.Range("F21").Value = 1984
$ws4.Range("G24").Value = "已售罄"
$ws4.Range("F28").Value = 632
$ws4.Range("F32").Value = 3475
$ws4.Range("F38").Value = 486
$ws4.Range("F42").Value = 319
$ws4.Range("F43").Value = 818
